# Applies the "test P7 with -10 percent" data update to the results workbook.
# Updates objValue/runtime/Z1 on 'general', assignment column on 'x',
# TBar values on 'TBar', and the Q cost table on 'Q'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 44.59727336721895
$ws.Range("B4").Value = 0.01399993896484375
$ws.Range("B6").Value = 44.59727336721895

$ws = $wb.Worksheets.Item("x")
$ws.Range("B3").Value = 3
$ws.Range("B5").Value = 8
$ws.Range("B6").Value = 7
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 12
$ws.Range("B10").Value = 11
$ws.Range("B13").Value = 6
$ws.Range("B14").Value = 13

$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B4").Value = 34.04101472405138
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 37.77398963608387
$ws.Range("B8").Value = 37.5860495735166
$ws.Range("B9").Value = 38.0993518219278
$ws.Range("B10").Value = 35.46173964959884
$ws.Range("B12").Value = 40.11094323173248
$ws.Range("B13").Value = 40.25111992906071
$ws.Range("B14").Value = 41.2899105680731
$ws.Range("B15").Value = 44.21573311673268

$ws = $wb.Worksheets.Item("Q")
$ws.Range("C12").Value = 128.3950000000002
$ws.Range("C13").Value = 116.7850000000001
$ws.Range("C14").Value = 119.6300000000002
$ws.Range("C15").Value = 117.9250000000001
$ws.Range("C16").Value = 126.7800000000001
$ws.Range("C17").Value = 46.91999999999942
$ws.Range("C18").Value = 36.10499999999942
$ws.Range("C19").Value = 34.91499999999942
$ws.Range("C20").Value = 37.48999999999942
$ws.Range("C21").Value = 39.43499999999941
$ws.Range("C22").Value = 250.5
$ws.Range("C23").Value = 276.8049999999997
$ws.Range("C24").Value = 255.4099999999997
$ws.Range("C25").Value = 266.8799999999997
$ws.Range("C26").Value = 255.55
$ws.Range("C32").Value = 262.7450000000006
$ws.Range("C33").Value = 276.5250000000006
$ws.Range("C34").Value = 245.9150000000006
$ws.Range("C35").Value = 272.2100000000005
$ws.Range("C36").Value = 252.4
$ws.Range("C37").Value = 211.470000000001
$ws.Range("C38").Value = 201.505000000001
$ws.Range("C39").Value = 197.6050000000011
$ws.Range("C40").Value = 207.6950000000011
$ws.Range("C41").Value = 208.635000000001
$ws.Range("C42").Value = 177.8700000000002
$ws.Range("C43").Value = 196.75
$ws.Range("C44").Value = 173.0900000000002
$ws.Range("C45").Value = 182.9250000000002
$ws.Range("C46").Value = 174.7350000000002
$ws.Range("C48").Value = 164.2850000000007
$ws.Range("C52").Value = 269.4200000000003
$ws.Range("C53").Value = 260.1850000000003
$ws.Range("C54").Value = 259.3350000000003
$ws.Range("C55").Value = 268.35
$ws.Range("C56").Value = 261.55
$ws.Range("C57").Value = 262.7450000000006
$ws.Range("C58").Value = 276.5250000000006
$ws.Range("C59").Value = 245.9150000000006
$ws.Range("C60").Value = 272.2100000000005
$ws.Range("C61").Value = 252.4
$ws.Range("C62").Value = 250.5
$ws.Range("C63").Value = 276.8049999999997
$ws.Range("C64").Value = 255.4099999999997
$ws.Range("C65").Value = 266.8799999999997
$ws.Range("C66").Value = 255.55
$ws.Range("C67").Value = 269.4200000000003
$ws.Range("C68").Value = 260.1850000000003
$ws.Range("C69").Value = 259.3350000000003
$ws.Range("C70").Value = 268.35
$ws.Range("C71").Value = 261.55
